$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = -20.46729999999997
$ws.Range("A12").Value = -22.49410000000004
$ws.Range("D13").Value = -7.9945
$ws.Range("A18").Value = -22.48150000000004
